$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.781.54'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '1.642.66'
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").Value = "'217.17"
$ws.Range("E5").Value = '  +1.09%  '

$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("E7").Value = '  +0.44%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  +0.12%  '

$ws.Range("D10").Value = "'19.16"
$ws.Range("E10").Value = '  +0.55%  '

$ws.Range("E11").Value = '  +0.07%  '

$ws.Range("D12").Value = '1.872.27'
$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").Value = '1.652.12'
$ws.Range("E13").Value = '  -0.16%  '

$ws.Range("E14").Value = '  -0.58%  '

$ws.Range("D16").Value = "'64.57"
$ws.Range("E16").Value = '  -0.64%  '

$ws.Range("D17").Value = '26.768.50'
$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("E18").Value = '  -1.49%  '

$ws.Range("D19").Value = "'214.11"
$ws.Range("E19").Value = '  -1.04%  '

$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("E21").Value = '  +0.59%  '

$ws.Range("D22").Value = "'2.35"
$ws.Range("E22").Value = '  +4.95%  '

$ws.Range("D23").Value = "'6.22"
$ws.Range("E23").Value = '  -0.47%  '

$ws.Range("D24").Value = "'9.28"
$ws.Range("E24").Value = '  -2.37%  '

$ws.Range("D25").Value = "'145.53"
$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("E26").Value = '  +0.43%  '

$ws.Range("E27").Value = '  -1.32%  '

$ws.Range("D28").Value = "'7.18"
$ws.Range("E28").Value = '  +0.24%  '

$ws.Range("D29").Value = "'15.61"
$ws.Range("E29").Value = '  -0.58%  '

$ws.Range("E30").Value = '  -1.54%  '

$ws.Range("E31").Value = '  +1.07%  '

$ws.Range("D32").Value = "'3.37"
$ws.Range("E32").Value = '  +0.12%  '

$ws.Range("D33").Value = "'3.00"
$ws.Range("E33").Value = '  -1.19%  '

$ws.Range("D34").Value = '1.283.71'
$ws.Range("E34").Value = '  +0.45%  '

$ws.Range("E35").Value = '  -0.11%  '

$ws.Range("E36").Value = '  +1.65%  '

$ws.Range("E37").Value = '  -0.64%  '

$ws.Range("E38").Value = '  +0.46%  '

$ws.Range("E39").Value = '  -1.33%  '

$ws.Range("E40").Value = '  +0.50%  '

$ws.Range("D41").Value = "'0.804"
$ws.Range("E41").Value = '  -1.13%  '

$ws.Range("E42").Value = '  -1.00%  '

$ws.Range("D43").Value = "'5.29"
$ws.Range("E43").Value = '  -2.68%  '

$ws.Range("D44").Value = '1.782.86'
$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").Value = "'61.33"
$ws.Range("E45").Value = '  +3.23%  '

$ws.Range("D46").Value = "'91.89"
$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").Value = "'1.60"
$ws.Range("E47").Value = '  +0.36%  '

$ws.Range("D48").Value = "'0.0517"
$ws.Range("E48").Value = '  +0.30%  '

$ws.Range("D49").Value = "'7.64"
$ws.Range("E49").Value = '  -1.89%  '

$ws.Range("D50").Value = "'0.0965"
$ws.Range("E50").Value = '  -0.07%  '

$ws.Range("E51").Value = '  +0.15%  '
